# Generate Report for Handback
# This script updates the localization-status workbook to reflect that the
# de-de and zh-cn handback packages have been generated / handed back.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

# --- Base GitHub URL used for the hyperlinks to the source markdown files ---
$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eaa460a0377430625298760da74a211f7a1c67f1/e2e/"
$mdFile1 = "01e95f41-cf93-4bb6-9df0-b5a2a251c996.md"
$mdFile2 = "7c333bfe-7b58-4c1d-a215-29ef7f507a30.md"

# ------------------------------------------------------------------
# 1. Overview sheet: status columns (zh-cn = E, de-de = F) for both rows
# ------------------------------------------------------------------
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ------------------------------------------------------------------
# 2. zh-cn sheet: Status, Latest Target File, Latest Handback File,
#    Latest Handback DateTime
# ------------------------------------------------------------------
$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("C3").Value = $statusText

$wsZhCn.Range("J2").Value = "01e95f41-cf93-4bb6-9df0-b5a2a251c996.9094f65ae79496ef1346984638288354f6f79385.zh-cn.xlf"
$wsZhCn.Range("J3").Value = "7c333bfe-7b58-4c1d-a215-29ef7f507a30.167675640e0b0aa216312997ae96db1176ced7af.zh-cn.xlf"

$wsZhCn.Range("K2").Value = "2016-10-18 13:41:16"
$wsZhCn.Range("K3").Value = "2016-10-18 13:41:16"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), ($baseUrl + $mdFile1), "", "", $mdFile1) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), ($baseUrl + $mdFile2), "", "", $mdFile2) | Out-Null

$zhCnI = $wsZhCn.Range("I2:I3")
$zhCnI.Font.Underline = 2
$zhCnI.Font.Color = 15570276

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ------------------------------------------------------------------
# 3. de-de sheet: Status, Latest Target File, Latest Handback File,
#    Latest Handback DateTime
# ------------------------------------------------------------------
$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("C3").Value = $statusText

$wsDeDe.Range("J2").Value = "01e95f41-cf93-4bb6-9df0-b5a2a251c996.9094f65ae79496ef1346984638288354f6f79385.de-de.xlf"
$wsDeDe.Range("J3").Value = "7c333bfe-7b58-4c1d-a215-29ef7f507a30.167675640e0b0aa216312997ae96db1176ced7af.de-de.xlf"

$wsDeDe.Range("K2").Value = "2016-10-18 13:41:34"
$wsDeDe.Range("K3").Value = "2016-10-18 13:41:34"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), ($baseUrl + $mdFile1), "", "", $mdFile1) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), ($baseUrl + $mdFile2), "", "", $mdFile2) | Out-Null

$deDeI = $wsDeDe.Range("I2:I3")
$deDeI.Font.Underline = 2
$deDeI.Font.Color = 15570276

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664
